# Daily attendance processing - 2025-11-28 21:24:13
# Move "System" to the front of the comma-separated "Recorded By" list
# in column G, preserving the relative order of the remaining entries.
# (Comparisons use .Equals() to stay case-sensitive, since this engine's
#  -eq/-contains/-ne operators are case-insensitive and would otherwise
#  also match the lowercase "system" token that appears in some rows.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2

    if ($val -and $val -like "*,*" -and $val -like "*System*") {
        $parts = $val -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        $hasSystem = $false
        foreach ($t in $trimmed) {
            if ($t.Equals("System")) {
                $hasSystem = $true
            }
        }

        if ($hasSystem -and -not $trimmed[0].Equals("System")) {
            $rest = @()
            $removed = $false
            foreach ($t in $trimmed) {
                if ((-not $removed) -and $t.Equals("System")) {
                    $removed = $true
                } else {
                    $rest += $t
                }
            }
            $new = @("System") + $rest
            $cell.Value2 = [string]::Join(", ", $new)
        }
    }
}
